$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 2.024603843688965
$ws.Range("B1").Value = 2.216437101364136
$ws.Range("C1").Value = 8.005729675292969
$ws.Range("D1").Value = 0.9627824425697327
$ws.Range("E1").Value = 0.5426732897758484
